$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental -> "false" (text, not boolean)
$ws.Range("B7").Formula = "=""false"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# Date -> updated timestamp
$ws.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# Case Sensitive -> "true" (text, not boolean)
$ws.Range("B15").Formula = "=""true"""
$ws.Range("B15").Copy()
$ws.Range("B15").PasteSpecial(-4163)
